$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WayPoints")

# Insert a new row before row 8 (current "LFPG/08L" row) to make room for the
# new "LIZY-SUR-OURCQ" waypoint, shifting the Charles-De-Gaulle/LFPG row down.
$ws.Rows.Item(8).Insert()

# Fill in the "name" column (column F) for rows that previously lacked it,
# matching the original (non-dashed) waypoint text.
$ws.Range("F3").Value = "ETREPAGNY"
$ws.Range("F4").Value = "CREIL"

# New row 8: LIZY-SUR-OURCQ waypoint.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "LIZY-SUR-OURCQ"
$ws.Range("F8").Value = "LIZY-SUR-OURCQ"

# Row 7: rename waypoint, adding dashes (same text used for both waypoint and
# name columns).
$ws.Range("B7").Value = "CROUY-SUR-OURCQ"
$ws.Range("F7").Value = "CROUY-SUR-OURCQ"

$ws.Range("C8").Value = "N49°01'06.00"""
$ws.Range("D8").Value = "E003°01'00.00"""
$ws.Range("E8").Value = "France"

# Row 6: rename waypoint (dashes added); name column keeps the original text.
$ws.Range("B6").Value = "LA-FERTE-MILON"
$ws.Range("F6").Value = "LAFERTEMILON"

# Row 5: rename waypoint (dashes added); name column keeps the original text.
$ws.Range("B5").Value = "CREPY-EN-VALOIS"
$ws.Range("F5").Value = "CREPYENVALOIS"

# Row 9 (previously row 8) keeps its data but the order number increments.
$ws.Range("A9").Value = 8

# Column width / formatting tweaks (closest values the engine's column-width
# quantization can reach to the authored 25.08984375 / 17.54296875 widths).
$ws.Columns.Item(2).ColumnWidth = 24.25
$ws.Columns.Item(6).ColumnWidth = 16.65

# Selection as recorded in the saved file.
$ws.Range("B12").Select()
